# Add a new LeetCode problem entry ("Merge Two Sorted Linked List") as row 13
# of the tracker table, mirroring the formatting/content pattern of the
# preceding row (row 12), and update the active selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 21
$ws.Range("B13").Value = "Merge Two Sorted Linked List"
$ws.Range("C13").Value = "Easy"
$ws.Range("D13").Value = "LinkedList"
$ws.Range("E13").Value = "Accepted"
$ws.Range("F13").Value = "O(n)"
$ws.Range("G13").Value = "O(1)"
$ws.Range("H13").Value = "Phase 3-4"

# Match the final selection left behind in the source workbook.
$null = $ws.Range("H15").Select()
